# data/gdp.xlsx: quarterly GDP value for 2024-06-30 (row 3) was corrected
# from 2.053 to 3.053; the % Change formula in C3 recalculates automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("B3").Value = 3.053
